$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Transmitance column (B) rows 3-18 from 1 to 100 (added the extra
# ordinary ray for polarimetric acquisitions)
for ($row = 3; $row -le 18; $row++) {
    $ws.Cells.Item($row, 2).Value = 100
}

# Update the active cell selection to B18
$ws.Range("B18").Select()
